# Update of existing "Mean" (G) / summer-period stats (L,M,N / I) figures
# for the 2016-2020, 2017-2021 and 2018-2022 reporting windows, plus append
# a brand-new 2019 - 2023 reporting block (rows 51-69) for
# "Kuku at N. Johnstone Farm Bridge".  (river update May 2024)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Tweak a handful of already-published cells -------------------------
# Each entry is: row, column index (1-based), new value
$cellEdits = @(
    @(3, 7, 0.0125889091403071),
    @(4, 7, 0.0125889091403071),
    @(6, 7, 0.014055513417007),
    @(6, 12, 0.00412),
    @(6, 13, 0.02232),
    @(6, 14, 0.04167),
    @(7, 7, 0.014055513417007),
    @(7, 12, 0.00412),
    @(7, 13, 0.02232),
    @(7, 14, 0.04167),
    @(18, 7, 0.0122239353294284),
    @(19, 7, 0.0122239353294284),
    @(21, 7, 0.0129260926087985),
    @(21, 13, 0.02011),
    @(21, 14, 0.03898),
    @(22, 7, 0.0129260926087985),
    @(22, 13, 0.02011),
    @(22, 14, 0.03898),
    @(33, 7, 0.0114680747220764),
    @(34, 7, 0.0114680747220764),
    @(40, 7, 0.0121505269851743),
    @(40, 9, 0.06156),
    @(40, 14, 0.03025),
    @(41, 7, 0.0121505269851743),
    @(41, 9, 0.06156),
    @(41, 14, 0.03025)
)

foreach ($edit in $cellEdits) {
    $ws.Cells.Item($edit[0], $edit[1]).Value = $edit[2]
}

# --- 2. Append the new 2019 - 2023 "RepSite" block (rows 51-69) ------------
$newRows = @(
    @("Kuku at N. Johnstone Farm Bridge", "ASPM", "D", "2019 - 2023", "RepSite", 0.157, 0.1544, 0.189, 0.189, "", "", 0.157, 0.17885, 0.189, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", ""),
    @("Kuku at N. Johnstone Farm Bridge", "DRP (95th Percentile)", "A", "2019 - 2023", "RepSite", 0.011, 0.0118035714285714, 0.027, 0.02, "", "", 0.0095, 0.01698, 0.01902, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "mg/L"),
    @("Kuku at N. Johnstone Farm Bridge", "DRP (Median)", "C", "2019 - 2023", "RepSite", 0.011, 0.0118035714285714, 0.027, 0.02, "", "", 0.0095, 0.01698, 0.01902, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "mg/L"),
    @("Kuku at N. Johnstone Farm Bridge", "E coli (>260)", "E", "2019 - 2023", "RepSite", 412, 849.196428571429, 7900, 4129.8, 32.1428571428571, 69.6428571428571, 380, 1196, 2515.32, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "% exceedances over 260/100 mL"),
    @("Kuku at N. Johnstone Farm Bridge", "E coli (>540)", "E", "2019 - 2023", "RepSite", 412, 849.196428571429, 7900, 4129.8, 32.1428571428571, 69.6428571428571, 380, 1196, 2515.32, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "% exceedances over 540/100 mL"),
    @("Kuku at N. Johnstone Farm Bridge", "E coli (Median)", "E", "2019 - 2023", "RepSite", 412, 849.196428571429, 7900, 4129.8, 32.1428571428571, 69.6428571428571, 380, 1196, 2515.32, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "E. coli/100 mL"),
    @("Kuku at N. Johnstone Farm Bridge", "E coli (95th Percentile)", "E", "2019 - 2023", "RepSite", 412, 849.196428571429, 7900, 4129.8, 32.1428571428571, 69.6428571428571, 380, 1196, 2515.32, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "E. coli/100 mL"),
    @("Kuku at N. Johnstone Farm Bridge", "MCI", "D", "2019 - 2023", "RepSite", 79, 80.026, 98.13, 98.13, "", "", 79, 91.7845, 98.13, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", ""),
    @("Kuku at N. Johnstone Farm Bridge", "Ammoniacal-N (95th Percentile)", "A", "2019 - 2023", "RepSite", 0.00618, 0.0071632453055254, 0.02620635060341, 0.01915, "", "", 0.00463, 0.01265, 0.01625, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "mg NH4-N/L"),
    @("Kuku at N. Johnstone Farm Bridge", "Ammoniacal-N (Median)", "A", "2019 - 2023", "RepSite", 0.00618, 0.0071632453055254, 0.02620635060341, 0.01915, "", "", 0.00463, 0.01265, 0.01625, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "mg NH4-N/L"),
    @("Kuku at N. Johnstone Farm Bridge", "Nitrate-N (95th Percentile)", "B", "2019 - 2023", "RepSite", 0.921, 1.075375, 2.63, 2.087, "", "", 1.175, 1.5182, 1.961, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "mg NO3-N/L"),
    @("Kuku at N. Johnstone Farm Bridge", "Nitrate-N (Median)", "A", "2019 - 2023", "RepSite", 0.921, 1.075375, 2.63, 2.087, "", "", 1.175, 1.5182, 1.961, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "mg NO3-N/L"),
    @("Kuku at N. Johnstone Farm Bridge", "QMCI", "D", "2019 - 2023", "RepSite", 4.162, 3.9004, 4.46, 4.46, "", "", 4.162, 4.453, 4.46, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", ""),
    @("Kuku at N. Johnstone Farm Bridge", "Soluble Inorganic Nitrogen (95th Percentile)", "", "2019 - 2023", "RepSite", 0.9425, 1.095625, 2.639, 2.1064, "", "", 1.191, 1.53224, 1.97508, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "g/m3"),
    @("Kuku at N. Johnstone Farm Bridge", "Soluble Inorganic Nitrogen (Median)", "", "2019 - 2023", "RepSite", 0.9425, 1.095625, 2.639, 2.1064, "", "", 1.191, 1.53224, 1.97508, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "g/m3"),
    @("Kuku at N. Johnstone Farm Bridge", "Total Nitrogen (95th Percentile)", "", "2019 - 2023", "RepSite", 1.16, 1.26232142857143, 2.83, 2.32, "", "", 1.325, 1.6292, 2.0736, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "g/m3"),
    @("Kuku at N. Johnstone Farm Bridge", "Total Nitrogen (Median)", "", "2019 - 2023", "RepSite", 1.16, 1.26232142857143, 2.83, 2.32, "", "", 1.325, 1.6292, 2.0736, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "g/m3"),
    @("Kuku at N. Johnstone Farm Bridge", "Total Phosphorus (95th Percentile)", "", "2019 - 2023", "RepSite", 0.03, 0.0344464285714286, 0.104, 0.062, "", "", 0.025, 0.046, 0.05808, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "g/m3"),
    @("Kuku at N. Johnstone Farm Bridge", "Total Phosphorus (Median)", "", "2019 - 2023", "RepSite", 0.03, 0.0344464285714286, 0.104, 0.062, "", "", 0.025, 0.046, 0.05808, 1785330, 5495215, "Horowhenua District", "Waiopehu", "Ohau", "Ohau_1b", "g/m3")
)

$targetRow = 51
foreach ($rowValues in $newRows) {
    for ($col = 1; $col -le 21; $col++) {
        $ws.Cells.Item($targetRow, $col).Value = $rowValues[$col - 1]
    }
    $targetRow = $targetRow + 1
}
